$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Michael%Irvine%NULL%1,         Daniel%Coombs%NULL%1,         Julianne%Skarha%NULL%1,         Brandon%del Pozo%NULL%1,         Josiah%Rich%NULL%1,         Faye%Taxman%NULL%1,         Traci C.%Green%Traci.c.green@gmail.com%1]'
$ws.Range("E3").Value = '[Shaun%Truelove%NULL%1,         Orit%Abrahim%NULL%2,         Orit%Abrahim%NULL%0,         Chiara%Altare%NULL%2,         Chiara%Altare%NULL%0,         Stephen A.%Lauer%NULL%2,         Stephen A.%Lauer%NULL%0,         Krya H.%Grantz%NULL%2,         Krya H.%Grantz%NULL%0,         Andrew S.%Azman%NULL%2,         Andrew S.%Azman%NULL%0,         Paul%Spiegel%NULL%2,         Paul%Spiegel%NULL%0,         Parveen%Parmar%NULL%3,         Parveen%Parmar%NULL%0,         Parveen%Parmar%NULL%0,         NULL%NULL%NULL%0,         NULL%NULL%NULL%0,         NULL%NULL%NULL%0]'
$ws.Range("E4").Value = '[ M.%Hariri%null%1,          H.% Rihawi%null%1,          S.% Safadi%null%1,          M. A.% McGlasson%null%1,          W. % Obaid%null%1]'
$ws.Range("E5").Value = '[ I.%Bojorquez%null%1,          C.% Infante%null%1,          I.% Vieitez%null%1,          S.% Larrea%null%1,          C. % Santoro%null%1]'
$ws.Range("E6").Value = '[Maren%Hintermeier%NULL%1,         Hande%Gencer%NULL%1,         Katja%Kajikhina%NULL%1,         Sven%Rohleder%NULL%1,         Claudia%Hövener%NULL%1,         Marie%Tallarek%NULL%1,         Jacob%Spallek%NULL%1,         Kayvan%Bozorgmehr%kayvan.bozorgmehr@uni-bielefeld.de%1]'
$ws.Range("E8").Value = '[MH%Chew%chew.min.hoe@singhealth.com.sg%1,         F.H.%Koh%NULL%1,         JT%Wu%NULL%1,         S.%Ngaserin%NULL%1,         A.%Ng%NULL%1,         BC%Ong%NULL%1,         V.J.%Lee%NULL%1]'
